# Updates the "苏州-漫展信息" workbook to the next generated snapshot:
#  - bump "想去人数" (F column) counters on the 展览 (Exhibition) sheet
#  - add the new "苏州·龙猫和他的朋友·动漫作品音乐会" concert row to the
#    (currently empty) 演出 (Performance) sheet
#  - mirror both of those changes on 全部类型 (All types), which is the
#    union of every category sorted by date, so the new concert row lands
#    between the 02-25 and 03-08 Suzhou exhibition rows there

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) 展览 (Exhibition) — counter bumps only, no structural changes
# ---------------------------------------------------------------------
$exhibition = $wb.Worksheets.Item("展览")

$fBumps = @{
    3  = 1111
    5  = 90
    6  = 56
    8  = 11373
    9  = 4318
    10 = 27
    11 = 29
    13 = 2521
    14 = 1076
    15 = 121
    16 = 24
    17 = 174
    18 = 496
    19 = 11271
    20 = 11140
    22 = 39
}

foreach ($row in $fBumps.Keys) {
    $exhibition.Cells.Item($row, 6).Value = $fBumps[$row]
}

# ---------------------------------------------------------------------
# 2) 演出 (Performance) — insert the brand-new concert row under the
#    header row
# ---------------------------------------------------------------------
$performance = $wb.Worksheets.Item("演出")
$performance.Rows.Item(2).Insert()

$pA = $performance.Cells.Item(2, 1)
$pA.Font.Bold = $true
$pA.HorizontalAlignment = -4108
$pA.VerticalAlignment = -4160
$pA.Borders.LineStyle = 1
$pA.Value = 1

$performance.Cells.Item(2, 2).NumberFormat = "@"
$performance.Cells.Item(2, 2).Value = "2024-03-03"
$performance.Cells.Item(2, 3).Value = "苏州·龙猫和他的朋友·动漫作品音乐会"
$performance.Cells.Item(2, 4).Value = "星湖街555号高教区(体育馆南侧) 苏州独墅湖影剧院"
$performance.Cells.Item(2, 5).Value = "2024.03.03 19:30-03.03 21:00"
$performance.Cells.Item(2, 6).Value = 0
$performance.Cells.Item(2, 7).Value = 60
$performance.Cells.Item(2, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81799"
$performance.Cells.Item(2, 9).Value = "//i0.hdslb.com/bfs/openplatform/202402/gqnOEjvJ1707214629948.jpeg"

# ---------------------------------------------------------------------
# 3) 全部类型 (All types) — same counter bumps as 展览, plus the same new
#    concert row inserted in date order (row 14, ahead of the 03-08
#    cosplay show), which pushes the former rows 14-26 down to 15-27
# ---------------------------------------------------------------------
$allTypes = $wb.Worksheets.Item("全部类型")

foreach ($row in $fBumps.Keys) {
    $allTypes.Cells.Item($row, 6).Value = $fBumps[$row]
}

$allTypes.Rows.Item(14).Insert()

$aA = $allTypes.Cells.Item(14, 1)
$aA.Font.Bold = $true
$aA.HorizontalAlignment = -4108
$aA.VerticalAlignment = -4160
$aA.Borders.LineStyle = 1
$aA.Value = 13

$allTypes.Cells.Item(14, 2).NumberFormat = "@"
$allTypes.Cells.Item(14, 2).Value = "2024-03-03"
$allTypes.Cells.Item(14, 3).Value = "苏州·龙猫和他的朋友·动漫作品音乐会"
$allTypes.Cells.Item(14, 4).Value = "星湖街555号高教区(体育馆南侧) 苏州独墅湖影剧院"
$allTypes.Cells.Item(14, 5).Value = "2024.03.03 19:30-03.03 21:00"
$allTypes.Cells.Item(14, 6).Value = 0
$allTypes.Cells.Item(14, 7).Value = 60
$allTypes.Cells.Item(14, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81799"
$allTypes.Cells.Item(14, 9).Value = "//i0.hdslb.com/bfs/openplatform/202402/gqnOEjvJ1707214629948.jpeg"

# The rows that used to be 14-26 (index column A = 13-25) shifted down to
# 15-27 when the new row was inserted above them; Excel keeps their old
# literal index values when it shifts cells, but this sheet's index
# column is always a plain "row - 1" sequence, so renumber them back in
# line.
for ($row = 15; $row -le 27; $row++) {
    $allTypes.Cells.Item($row, 1).Value = $row - 1
}
